# Update the cached "Date and time (automatic)" footer field from
# 6/1/2023 -> 6/12/2023 everywhere it appears: the slide master and
# every slide layout's Date Placeholder.

$p = $ppt.ActivePresentation

$oldDate = "6/1/2023"
$newDate = "6/12/2023"

function Update-DatePlaceholder($shapes) {
    foreach ($sh in $shapes) {
        if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master's own Date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout under the master has its own Date placeholder copy.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}
